# Add the new "RegistrationPage" worksheet at the end of the workbook (after "HomePage"),
# populate it with the registration form data, and make it the active sheet/tab -
# mirroring a user adding a new sheet in the Excel UI.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "RegistrationPage"

# Column A is a bit wider than the default, like the other sheets in this workbook.
$ws.Columns.Item(1).ColumnWidth = 15.5

# Fill in column A first (top to bottom), then column B - this mirrors the order the
# shared-string table was populated in when the sheet was originally authored.
$ws.Range("A1").Value = "Gender"
$ws.Range("B1").Value = "Female"
$ws.Range("A2").Value = "First name"
$ws.Range("B2").Value = "First"
$ws.Range("A3").Value = "Last name"
$ws.Range("B3").Value = "Last"
$ws.Range("A4").Value = "Password"
$ws.Range("A5").Value = "Confirm password"
$ws.Range("B4").Value = "test123"
$ws.Range("B5").Value = "test123"

$ws.PageSetup.Orientation = 1

# Match the recorded selection / active cell for the new sheet, and make it the
# active tab (adding it after the last sheet + selecting it does this already).
$ws.Range("G11").Select()
